$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay as text (matches source inlineStr type)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.745.93"
$ws.Range("E2").Value = "  -0.39%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.294.36"
$ws.Range("E3").Value = "  -0.30%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.34"
$ws.Range("E5").Value = "  +1.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.05"
$ws.Range("E6").Value = "  -1.68%  "

$ws.Range("E7").Value = "  -1.94%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E9").Value = "  -1.91%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.91"
$ws.Range("E10").Value = "  -2.54%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0780"
$ws.Range("E11").Value = "  -0.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "18.68"
$ws.Range("E12").Value = "  +5.74%  "

$ws.Range("E13").Value = "  +2.29%  "

$ws.Range("E14").Value = "  +0.86%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.651.72"
$ws.Range("E15").Value = "  -0.28%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.298.70"
$ws.Range("E16").Value = "  -0.29%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.774"
$ws.Range("E17").Value = "  -0.55%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.650.06"
$ws.Range("E18").Value = "  -0.53%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.69"
$ws.Range("E19").Value = "  +0.94%  "

$ws.Range("E20").Value = "  -1.82%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.03"
$ws.Range("E22").Value = "  -1.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.12"
$ws.Range("E23").Value = "  -2.45%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.14"
$ws.Range("E24").Value = "  +0.08%  "

$ws.Range("E25").Value = "  +0.24%  "

$ws.Range("E26").Value = "  -1.50%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.67"
$ws.Range("E27").Value = "  -1.66%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "167.20"
$ws.Range("E28").Value = "  +0.75%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.06"
$ws.Range("E29").Value = "  +0.97%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.99"
$ws.Range("E30").Value = "  -0.48%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.70"
$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("E32").Value = "  +0.05%  "

$ws.Range("E33").Value = "  +1.76%  "

$ws.Range("E34").Value = "  -0.92%  "

$ws.Range("E35").Value = "  -6.60%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.35"
$ws.Range("E36").Value = "  -1.64%  "

$ws.Range("E37").Value = "  -0.49%  "

$ws.Range("E38").Value = "  -0.58%  "

$ws.Range("E39").Value = "  -1.47%  "

$ws.Range("E40").Value = "  -1.04%  "

$ws.Range("E41").Value = "  -2.16%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.992.73"
$ws.Range("E42").Value = "  -0.66%  "

$ws.Range("E43").Value = "  -1.92%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.19"
$ws.Range("E44").Value = "  +0.42%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.20"
$ws.Range("E45").Value = "  +6.11%  "

$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.13"
$ws.Range("E46").Value = "  -0.91%  "

$ws.Range("E47").Value = "  -0.63%  "

$ws.Range("B48").Value = "MultiversX"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.43"
$ws.Range("E48").Value = "  +0.23%  "

$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.84"
$ws.Range("E49").Value = "  +4.88%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.519.07"
$ws.Range("E50").Value = "  -0.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.76"
$ws.Range("E51").Value = "  -1.84%  "
